# Auto-update predictions and index for 2025-10-16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
,    @("Wed Oct 15", "Jamaica ✓ - Bermuda: 4:0", 3.62, "Jamaica", 4.5, "71%", "✓", 4, $true)
    ,    @("Wed Oct 15", "FBC Melgar ✓ - Alianza Universidad: 2:1", 1.74, "FBC Melgar", 2.5, "70%", "✓", 3, $false)
    ,    @("Wed Oct 15", "Puerto Rico - Argentina ✓: 0:6", 1.49, "Argentina", 2.5, "68%", "✓", 6, $false)
    ,    @("Wed Oct 15", "Sociedade Esportiva Palmeiras ✓ - Red Bull Bragantino: 5:1", 2.04, "Sociedade Esportiva Palmeiras", 3.5, "68%", "✓", 6, $false)
    ,    @("Wed Oct 15", "Club Athletico Paranaense  - Avaí FC: 1:1", 1.81, "Club Athletico Paranaense", 2.5, "65%", $null, 2, $true)
    ,    @("Wed Oct 15", "Al-Nahda Club  - Sur SC: 14:10", 1, "Al-Nahda Club", 1.5, "61%", $null, 24, $false)
    ,    @("Wed Oct 15", "Atlético Nacional ✓ - Deportivo Cali: 2:1", 1.36, "Atlético Nacional", 2.5, "61%", "✓", 3, $false)
    ,    @("Wed Oct 15", "United States ✓ - Australia: 2:1", 1.34, "United States", 2.5, "58%", "✓", 3, $false)
    ,    @("Wed Oct 15", "Rodina Moscow X - FK Chelyabinsk: 4:5", 1.1, "Rodina Moscow", 2.5, "49%", "X", 9, $false)
    ,    @("Wed Oct 15", "Panama  - Suriname: 1:1", 1.38, "Panama", 2.5, "41%", $null, 2, $true)
    ,    @("Wed Oct 15", "Ironi Modiin - Maccabi Petah Tikva ✓: 1:2", 3.04, "Maccabi Petah Tikva", 4.5, "36%", "✓", 3, $true)
    ,    @("Wed Oct 15", "Shkendija Tetovo  - Tikves Kavadarci: 13:30", 0.98, "Shkendija Tetovo", 1.5, "33%", $null, 43, $false)
    ,    @("Wed Oct 15", "Club Sport Sébaco  - ART Municipal Jalapa: 2:2", 1.43, "Club Sport Sébaco", 2.5, "32%", $null, 4, $false)
    ,    @("Wed Oct 15", "Al-Shabab SC (Oman)  - Al-Nasr Club: 14:00", 1.72, "Al-Shabab SC (Oman)", 2.5, "31%", $null, 14, $false)
    ,    @("Wed Oct 15", "22 de Julio FC X - Cumbayá FC: 1:3", 1.44, "22 de Julio FC", 2.5, "29%", "X", 4, $false)
    ,    @("Wed Oct 15", "Mexico  - Ecuador: 1:1", 0, "Mexico", 0.5, "24%", $null, 2, $false)
    ,    @("Wed Oct 15", "Curaçao  - Trinidad and Tobago: 1:1", 2.24, "Curaçao", 3.5, "24%", $null, 2, $true)
    ,    @("Wed Oct 15", "USM Khenchela - MC Algiers : 23:00", 1.13, "MC Algiers", 2.5, "23%", $null, 23, $false)
    ,    @("Wed Oct 15", "Canada - Colombia : 0:0", 2.97, "Colombia", 3.5, "23%", $null, 0, $true)
    ,    @("Wed Oct 15", "Belize - Venezuela : 00:00", 2.45, "Venezuela", 3.5, "21%", $null, 0, $true)
    ,    @("Wed Oct 15", "JS Saoura - JS Kabylie : 23:00", 1.42, "JS Kabylie", 2.5, "20%", $null, 23, $false)
    ,    @("Wed Oct 15", "Arsenal Tula ✓ - Sokol Saratov: 1:0", 0.71, "Arsenal Tula", 1.5, "19%", "✓", 1, $true)
    ,    @("Wed Oct 15", "AD Municipal Liberia  - Puntarenas FC: 0:0", 1.24, "AD Municipal Liberia", 2.5, "18%", $null, 0, $true)
    ,    @("Wed Oct 15", "ASO Chlef - CR Belouizdad : 23:00", 1.4, "CR Belouizdad", 2.5, "18%", $null, 23, $false)
    ,    @("Wed Oct 15", "Botafogo de Futebol e Regatas - CR Flamengo ✓: 0:3", 1.43, "CR Flamengo", 2.5, "18%", "✓", 3, $false)
    ,    @("Wed Oct 15", "MFK Zvolen ✓ - Dynamo Malzenice: 3:1", 2.02, "MFK Zvolen", 3.5, "16%", "✓", 4, $false)
    ,    @("Wed Oct 15", "MB Rouissat  - ES Ben Aknoun: 23:00", 1.27, "MB Rouissat", 2.5, "16%", $null, 23, $false)
    ,    @("Wed Oct 15", "CS Constantine  - Olympique Akbou: 23:00", 1.14, "CS Constantine", 2.5, "15%", $null, 23, $false)
    ,    @("Wed Oct 15", "Uganda People's Defence Force FC X - Maroons FC: 1:2", 2.06, "Uganda People's Defence Force FC", 3.5, "14%", "X", 3, $true)
    ,    @("Wed Oct 15", "Shinnik Yaroslavl X - Neftekhimik Nizhnekamsk: 4:6", 0.86, "Shinnik Yaroslavl", 1.5, "14%", "X", 10, $false)
    ,    @("Wed Oct 15", "FK Slutsk - Neman Grodno X: 2:1", 1.16, "Neman Grodno", 2.5, "10%", "X", 3, $false)
    ,    @("Wed Oct 15", "Dhofar Club  - Samail SC: 15:25", 0.66, "Dhofar Club", 1.5, "10%", $null, 40, $false)
    ,    @("Wed Oct 15", "Chaika Peschanokopskoe X - KamAZ Naberezhnye Chelny: 1:3", 1.54, "Chaika Peschanokopskoe", 2.5, "10%", "X", 4, $false)
    ,    @("Wed Oct 15", "MC El Bayadh - MC Oran : 23:00", 1.1, "MC Oran", 2.5, "9%", $null, 23, $false)
    ,    @("Wed Oct 15", "El Salvador - Guatemala ✓: 0:1", 0.94, "Guatemala", 1.5, "6%", "✓", 1, $true)
    ,    @("Wed Oct 15", "Al-Rustaq Club - Saham Club : 16:40", 1.18, "Saham Club", 2.5, "6%", $null, 56, $false)
    ,    @("Wed Oct 15", "Paradou AC  - ES Mostaganem: 23:00", 1.65, "Paradou AC", 2.5, "6%", $null, 23, $false)
    ,    @("Wed Oct 15", "CD Choloma - Lobos de la UPNFM X: 5:1", 2.52, "Lobos de la UPNFM", 3.5, "5%", "X", 6, $false)
    ,    @("Wed Oct 15", "ES Sétif - USM Alger : 23:00", 0.71, "USM Alger", 1.5, "4%", $null, 23, $false)
    ,    @("Wed Oct 15", "Platense FC ✓ - CD Victoria La Ceiba: 2:0", 1.9, "Platense FC", 2.5, "3%", "✓", 2, $true)
)

$startRow = 2
$r = $startRow
foreach ($row in $rows) {
    for ($c = 1; $c -le 9; $c++) {
        $val = $row[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($val -eq $null) {
            $cell.Value = ""
        } elseif ($val -is [string]) {
            # leading apostrophe forces literal text so look-alike
            # percentages / times are not auto-converted to numbers
            $cell.Value = "'" + $val
        } else {
            $cell.Value = $val
        }
    }
    $r++
}

